$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

# Columns A and C look like a date / a number respectively, so Excel would
# otherwise auto-convert them. Force text storage (matching the source's
# t="str" cells) by briefly switching to a text format, then clear the
# format again so the cell keeps the default style (no explicit "s" index),
# just like the rest of the table.
$cA = $ws.Cells.Item($row, 1)
$cA.NumberFormat = "@"
$cA.Value = "2025-10-15"
$cA.ClearFormats()

$ws.Cells.Item($row, 2).Value = "Pick 3"

$cC = $ws.Cells.Item($row, 3)
$cC.NumberFormat = "@"
$cC.Value = "251015"
$cC.ClearFormats()

$ws.Cells.Item($row, 4).Value = "0-2-5"
$ws.Cells.Item($row, 5).Value = "2025-10-15T21:37:48.664+04:00"
